$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 2.63
$ws.Range("K2").Value = 1.91

# Row 3
$ws.Range("BC3").Value = 126

# Row 5
$ws.Range("H5").Value = 5
$ws.Range("N5").Value = 13
$ws.Range("O5").Value = 1.22
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 1.75
$ws.Range("R5").Value = 2.05
$ws.Range("W5").Value = 6.5
$ws.Range("Y5").Value = 9
$ws.Range("Z5").Value = 8
$ws.Range("AC5").Value = 11
$ws.Range("AD5").Value = 9.5
$ws.Range("AE5").Value = 23
$ws.Range("AG5").Value = 21
$ws.Range("AH5").Value = 41
$ws.Range("AI5").Value = 29
$ws.Range("AJ5").Value = 126
$ws.Range("AK5").Value = 67
$ws.Range("AL5").Value = 67
$ws.Range("AO5").Value = 6
$ws.Range("AQ5").Value = 17
$ws.Range("AU5").Value = 10
$ws.Range("AV5").Value = 67
$ws.Range("AW5").Value = 10
$ws.Range("AX5").Value = 41
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 251
$ws.Range("BA5").Value = 251

# Row 8
$ws.Range("N8").Value = 9
$ws.Range("Q8").Value = 2.15
$ws.Range("R8").Value = 1.67

# Row 11
$ws.Range("N11").Value = 9
